$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "abbreviation" column between "name" and "age", with
# single-letter abbreviations for each greek-letter row, and renumber
# the ages to 1/2/3.

$ws.Range("C1").Value = $ws.Range("B1").Value2
$ws.Range("B1").Value = "abbreviation"

$ws.Range("C2").Value = $ws.Range("B2").Value2
$ws.Range("C3").Value = $ws.Range("B3").Value2
$ws.Range("C4").Value = $ws.Range("B4").Value2

$ws.Range("B2").Value = "a"
$ws.Range("B3").Value = "b"
$ws.Range("B4").Value = "g"

$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 3

$ws.Range("C16").Select()
